# Commit: theme swap (Integral -> Office on the slide master's theme) and
# updated table-style IDs on the three summary tables (slides 14-16).
#
# The deck's slide master / all slide layouts / the package-level theme
# relationship all resolve to ppt/theme/theme1.xml, which originally held
# the "Integral" (Red Violet) theme. The author re-applied the built-in
# "Office" design to the deck, which replaced the color scheme used by
# every slide (the font/format schemes were already identical between the
# two themes, so only the 10 non-black/white theme colors actually change).

$p = $ppt.ActivePresentation

# --- 1. Re-point the slide master's theme colors at the "Office" palette ---
# Any slide's ThemeColorScheme writes straight through to the shared
# ppt/theme/theme1.xml color scheme used by the whole deck, in the fixed
# index order dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}

# --- 2. Swap the custom table style for the built-in one on the three
#        balance-sheet tables (slides 14, 15 and 16) ---
$newTableStyleId = "{01E3AA2A-B34F-496A-963C-0BD6CF6E1852}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}
